$wb = $excel.ActiveWorkbook

# Add the new "Translations question" sheet right after "Translations" by
# copying the existing sheet (this keeps the column 1 bestFit width/format
# identical to the original, matching the real author's workflow) and then
# overwrite its contents for the new, smaller translation-export sample.
$ws1 = $wb.Worksheets.Item(1)
$ws1.Copy($null, $ws1)
$newWs = $wb.Worksheets.Item(2)
$newWs.Name = "Translations question"

# Drop the extra sample rows (3-5) copied from "Translations" - the new
# sheet only needs a header row plus a single data row.
$newWs.Rows.Item(3).Resize(3).Delete()

# Overwrite the header + single data row with the combobox-option sample.
$newWs.Range("A1").Value = "Entity Id"
$newWs.Range("B1").Value = "Type"
$newWs.Range("C1").Value = "Index"
$newWs.Range("D1").Value = "Original"
$newWs.Range("E1").Value = "Translation"

# New shared strings must be interned in this order: OptionTitle (already
# exists), Combobox Option, Опция Комбобокса, then the 1111...1 entity id.
$newWs.Range("B2").Value = "OptionTitle"
$newWs.Range("C2").Value = 1
$newWs.Range("D2").Value = "Combobox Option"
$newWs.Range("E2").Value = "Опция Комбобокса"

# A2 needs Text number format since the value looks numeric-like.
$newWs.Range("A2").NumberFormat = "@"
$newWs.Range("A2").Value = "11111111111111111111111111111111"

# Approximate the new sheet's resized columns - B keeps bestFit content but
# is far narrower now ("OptionTitle"/"Type" vs the old "ValidationMessage"),
# and C/D/E were dragged to custom widths in the authored workbook.
$newWs.Columns.Item(2).ColumnWidth = 10.26
$newWs.Columns.Item(3).ColumnWidth = 10.67
$newWs.Columns.Item(4).ColumnWidth = 16
$newWs.Columns.Item(5).ColumnWidth = 17.67

# Selection moves to A3 on the new, now-active sheet.
$newWs.Range("A3").Select()
